$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update description text in A2
$ws.Range("A2").Value = "Results are input for bridge_plain.png.csv and other conductivity mapping  CSV files."

# Swap the Lambda_i / Lambda_e labels (F9 <-> F10)
$ws.Range("F9").Value = "Lambda_e"
$ws.Range("F10").Value = "Lambda_i"

# Update the selection / active cell to A3
$ws.Range("A3").Select()
